# Trade #16 closed at 2026-02-16 22:53:38 - base_strategy DOWN +0.000%
# Appends a new trade row (row 17) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the existing OPEN-trade row layout.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column A: Trade #
    $ws.Range("A17").Value = 16

    # Column B: Date - force text so Excel doesn't auto-convert to a date serial
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("B17").Value = "2026-02-16"
    $ws.Range("B17").Style = "Normal"

    # Column C: Time - force text so Excel doesn't auto-convert to a time serial
    $ws.Range("C17").NumberFormat = "@"
    $ws.Range("C17").Value = "22:53:38"
    $ws.Range("C17").Style = "Normal"

    # Column D: Strategy
    $ws.Range("D17").Value = "base_strategy"

    # Column E: Side
    $ws.Range("E17").Value = "DOWN"

    # Column F: Entry Price
    $ws.Range("F17").Value = 49.999998

    # Column G: Exit Price - blank (trade still open), keep as an empty text cell
    $ws.Range("G17").Value = "'"
    $ws.Range("G17").Style = "Normal"

    # Column H: Status
    $ws.Range("H17").Value = "OPEN"

    # Column I: P&L %
    $ws.Range("I17").Value = 0

    # Column J: P&L $
    $ws.Range("J17").Value = 0

    # Column K: Capital After
    $ws.Range("K17").Value = 100

    # Column L: Entry Slippage (bps)
    $ws.Range("L17").Value = 0

    # Column M: Exit Slippage (bps)
    $ws.Range("M17").Value = 0

    # Column N: Confidence
    $ws.Range("N17").Value = 0.6

    # Column O: Entry Reason
    $ws.Range("O17").Value = "Normal spread capture: 19600 bps"

    # Column P: Exit Reason - blank (trade still open), keep as an empty text cell
    $ws.Range("P17").Value = "'"
    $ws.Range("P17").Style = "Normal"

    # Column Q: Duration (min)
    $ws.Range("Q17").Value = 0
}
